# Split utils.py into several files: add a new "utils" leaf row under the
# "imylu" folder, and update / reword several of the existing description
# comments. Mirrors docs_cn/dictionary_tree.xlsx commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row for "utils" just above the current row 14 (the
#        "pic" row), shifting it (and everything below) down by one. We only
#        touch columns B:F so the used range doesn't balloon out.
$ws.Range("B14:F14").Insert(-4121) # xlShiftDown

# Give the freshly inserted row 14 the same per-cell formatting as row 13
# (the row immediately above it, which has the same "leaf under imylu" style).
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = $ws.Rows.Item(13).RowHeight
$excel.CutCopyMode = 0

# --- 2. Fill in the new "utils" row.
$ws.Range("D14").Value = "utils"
$ws.Range("E14").Value = " >>> Utils functions and classes."

# --- 3. Reword the comments (column E) for the existing algorithm folders
#        now that utils.py has been split up.
$ws.Range("E5").Value = " >>> K-Means."
$ws.Range("E8").Value = " >>> GBDT, Random Forest and Isolation Forest."
$ws.Range("E9").Value = " >>> Linear Regression, Logistic Regression and Ridge. "
$ws.Range("E10").Value = " >>> KNN. "
$ws.Range("E11").Value = " >>> HMM and Gaussian Naive Bayes."
$ws.Range("E12").Value = " >>> ALS. "
$ws.Range("E13").Value = " >>> Decision Tree, Regression Tree and Isolation Tree."

# --- 4. Update the sheet's last-used selection to reflect where the author
#        ended up after editing.
$ws.Range("E18").Select() | Out-Null
